$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 'ECs'
$ws.Cells.Item(2,2).Value = 'Dnajb11'
$ws.Cells.Item(2,3).Value = 'Prtg'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 8.465110333333334
$ws.Cells.Item(2,8).Value = 25.395331
$ws.Cells.Item(2,9).Value = 0.2119143389236239
$ws.Cells.Item(2,10).Value = 0.2119143389236239
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.1416426666666667
$ws.Cells.Item(2,14).Value = 0.424928
$ws.Cells.Item(2,15).Value = 0.06389210658948997
$ws.Cells.Item(2,16).Value = 0.06389210658948997
$ws.Cells.Item(2,17).Value = 1.199020801240889
$ws.Cells.Item(2,18).Value = 10.791187211168
$ws.Cells.Item(2,19).Value = 0.01353965353034948
$ws.Cells.Item(2,20).Value = 0.01353965353034948

$ws.Cells.Item(3,1).Value = 'ECs'
$ws.Cells.Item(3,2).Value = 'Dnajb11'
$ws.Cells.Item(3,3).Value = 'Prtg'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 8.465110333333334
$ws.Cells.Item(3,8).Value = 25.395331
$ws.Cells.Item(3,9).Value = 0.2119143389236239
$ws.Cells.Item(3,10).Value = 0.2119143389236239
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.357753333333333
$ws.Cells.Item(3,14).Value = 4.073259999999999
$ws.Cells.Item(3,15).Value = 0.6124547266518232
$ws.Cells.Item(3,16).Value = 0.6124547266518232
$ws.Cells.Item(3,17).Value = 11.49353177211778
$ws.Cells.Item(3,18).Value = 103.44178594906
$ws.Cells.Item(3,19).Value = 0.1297879385190699
$ws.Cells.Item(3,20).Value = 0.1297879385190699

$ws.Cells.Item(4,1).Value = 'ECs'
$ws.Cells.Item(4,2).Value = 'Dnajb11'
$ws.Cells.Item(4,3).Value = 'Prtg'
$ws.Cells.Item(4,4).Value = 'MuSCs'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 8.465110333333334
$ws.Cells.Item(4,8).Value = 25.395331
$ws.Cells.Item(4,9).Value = 0.2119143389236239
$ws.Cells.Item(4,10).Value = 0.2119143389236239
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.7173469999999998
$ws.Cells.Item(4,14).Value = 2.152041
$ws.Cells.Item(4,15).Value = 0.3235805429554008
$ws.Cells.Item(4,16).Value = 0.3235805429554009
$ws.Cells.Item(4,17).Value = 6.072421502285666
$ws.Cells.Item(4,18).Value = 54.65179352057099
$ws.Cells.Item(4,19).Value = 0.06857135684894104
$ws.Cells.Item(4,20).Value = 0.06857135684894106

$ws.Cells.Item(5,1).Value = 'ECs'
$ws.Cells.Item(5,2).Value = 'Dnajb11'
$ws.Cells.Item(5,3).Value = 'Prtg'
$ws.Cells.Item(5,4).Value = 'Resolving-Mac'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 8.465110333333334
$ws.Cells.Item(5,8).Value = 25.395331
$ws.Cells.Item(5,9).Value = 0.2119143389236239
$ws.Cells.Item(5,10).Value = 0.2119143389236239
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.000161
$ws.Cells.Item(5,14).Value = 0.000483
$ws.Cells.Item(5,15).Value = 0.00007262380328602412
$ws.Cells.Item(5,16).Value = 0.00007262380328602412
$ws.Cells.Item(5,17).Value = 0.001362882763666667
$ws.Cells.Item(5,18).Value = 0.012265944873
$ws.Cells.Item(5,19).Value = 0.0000153900252634771
$ws.Cells.Item(5,20).Value = 0.0000153900252634771

$ws.Cells.Item(6,1).Value = 'FAPs'
$ws.Cells.Item(6,2).Value = 'Dnajb11'
$ws.Cells.Item(6,3).Value = 'Prtg'
$ws.Cells.Item(6,4).Value = 'ECs'
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 13.040437
$ws.Cells.Item(6,8).Value = 39.121311
$ws.Cells.Item(6,9).Value = 0.3264524001829507
$ws.Cells.Item(6,10).Value = 0.3264524001829507
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.1416426666666667
$ws.Cells.Item(6,14).Value = 0.424928
$ws.Cells.Item(6,15).Value = 0.06389210658948997
$ws.Cells.Item(6,16).Value = 0.06389210658948997
$ws.Cells.Item(6,17).Value = 1.847082271178667
$ws.Cells.Item(6,18).Value = 16.623740440608
$ws.Cells.Item(6,19).Value = 0.02085773154888393
$ws.Cells.Item(6,20).Value = 0.02085773154888393

$ws.Cells.Item(7,1).Value = 'FAPs'
$ws.Cells.Item(7,2).Value = 'Dnajb11'
$ws.Cells.Item(7,3).Value = 'Prtg'
$ws.Cells.Item(7,4).Value = 'FAPs'
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 13.040437
$ws.Cells.Item(7,8).Value = 39.121311
$ws.Cells.Item(7,9).Value = 0.3264524001829507
$ws.Cells.Item(7,10).Value = 0.3264524001829507
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.357753333333333
$ws.Cells.Item(7,14).Value = 4.073259999999999
$ws.Cells.Item(7,15).Value = 0.6124547266518232
$ws.Cells.Item(7,16).Value = 0.6124547266518232
$ws.Cells.Item(7,17).Value = 17.70569680487333
$ws.Cells.Item(7,18).Value = 159.35127124386
$ws.Cells.Item(7,19).Value = 0.1999373155188807
$ws.Cells.Item(7,20).Value = 0.1999373155188807

$ws.Cells.Item(8,1).Value = 'FAPs'
$ws.Cells.Item(8,2).Value = 'Dnajb11'
$ws.Cells.Item(8,3).Value = 'Prtg'
$ws.Cells.Item(8,4).Value = 'MuSCs'
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 13.040437
$ws.Cells.Item(8,8).Value = 39.121311
$ws.Cells.Item(8,9).Value = 0.3264524001829507
$ws.Cells.Item(8,10).Value = 0.3264524001829507
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.7173469999999998
$ws.Cells.Item(8,14).Value = 2.152041
$ws.Cells.Item(8,15).Value = 0.3235805429554008
$ws.Cells.Item(8,16).Value = 0.3235805429554009
$ws.Cells.Item(8,17).Value = 9.354518360638997
$ws.Cells.Item(8,18).Value = 84.19066524575098
$ws.Cells.Item(8,19).Value = 0.105633644900293
$ws.Cells.Item(8,20).Value = 0.105633644900293

$ws.Cells.Item(9,1).Value = 'FAPs'
$ws.Cells.Item(9,2).Value = 'Dnajb11'
$ws.Cells.Item(9,3).Value = 'Prtg'
$ws.Cells.Item(9,4).Value = 'Resolving-Mac'
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 13.040437
$ws.Cells.Item(9,8).Value = 39.121311
$ws.Cells.Item(9,9).Value = 0.3264524001829507
$ws.Cells.Item(9,10).Value = 0.3264524001829507
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.000161
$ws.Cells.Item(9,14).Value = 0.000483
$ws.Cells.Item(9,15).Value = 0.00007262380328602412
$ws.Cells.Item(9,16).Value = 0.00007262380328602412
$ws.Cells.Item(9,17).Value = 0.002099510357
$ws.Cells.Item(9,18).Value = 0.018895593213
$ws.Cells.Item(9,19).Value = 0.00002370821489313704
$ws.Cells.Item(9,20).Value = 0.00002370821489313704

$ws.Cells.Item(10,1).Value = 'MuSCs'
$ws.Cells.Item(10,2).Value = 'Dnajb11'
$ws.Cells.Item(10,3).Value = 'Prtg'
$ws.Cells.Item(10,4).Value = 'ECs'
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 7.170964000000001
$ws.Cells.Item(10,8).Value = 21.512892
$ws.Cells.Item(10,9).Value = 0.1795168681406561
$ws.Cells.Item(10,10).Value = 0.1795168681406561
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.1416426666666667
$ws.Cells.Item(10,14).Value = 0.424928
$ws.Cells.Item(10,15).Value = 0.06389210658948997
$ws.Cells.Item(10,16).Value = 0.06389210658948997
$ws.Cells.Item(10,17).Value = 1.015714463530667
$ws.Cells.Item(10,18).Value = 9.141430171775999
$ws.Cells.Item(10,19).Value = 0.01146971087385422
$ws.Cells.Item(10,20).Value = 0.01146971087385422

$ws.Cells.Item(11,1).Value = 'MuSCs'
$ws.Cells.Item(11,2).Value = 'Dnajb11'
$ws.Cells.Item(11,3).Value = 'Prtg'
$ws.Cells.Item(11,4).Value = 'FAPs'
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 7.170964000000001
$ws.Cells.Item(11,8).Value = 21.512892
$ws.Cells.Item(11,9).Value = 0.1795168681406561
$ws.Cells.Item(11,10).Value = 0.1795168681406561
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.357753333333333
$ws.Cells.Item(11,14).Value = 4.073259999999999
$ws.Cells.Item(11,15).Value = 0.6124547266518232
$ws.Cells.Item(11,16).Value = 0.6124547266518232
$ws.Cells.Item(11,17).Value = 9.736400274213333
$ws.Cells.Item(11,18).Value = 87.62760246791998
$ws.Cells.Item(11,19).Value = 0.1099459544064769
$ws.Cells.Item(11,20).Value = 0.1099459544064769

$ws.Cells.Item(12,1).Value = 'MuSCs'
$ws.Cells.Item(12,2).Value = 'Dnajb11'
$ws.Cells.Item(12,3).Value = 'Prtg'
$ws.Cells.Item(12,4).Value = 'MuSCs'
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 7.170964000000001
$ws.Cells.Item(12,8).Value = 21.512892
$ws.Cells.Item(12,9).Value = 0.1795168681406561
$ws.Cells.Item(12,10).Value = 0.1795168681406561
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.7173469999999998
$ws.Cells.Item(12,14).Value = 2.152041
$ws.Cells.Item(12,15).Value = 0.3235805429554008
$ws.Cells.Item(12,16).Value = 0.3235805429554009
$ws.Cells.Item(12,17).Value = 5.144069512508
$ws.Cells.Item(12,18).Value = 46.29662561257199
$ws.Cells.Item(12,19).Value = 0.05808816566260661
$ws.Cells.Item(12,20).Value = 0.05808816566260662

$ws.Cells.Item(13,1).Value = 'MuSCs'
$ws.Cells.Item(13,2).Value = 'Dnajb11'
$ws.Cells.Item(13,3).Value = 'Prtg'
$ws.Cells.Item(13,4).Value = 'Resolving-Mac'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 7.170964000000001
$ws.Cells.Item(13,8).Value = 21.512892
$ws.Cells.Item(13,9).Value = 0.1795168681406561
$ws.Cells.Item(13,10).Value = 0.1795168681406561
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.000161
$ws.Cells.Item(13,14).Value = 0.000483
$ws.Cells.Item(13,15).Value = 0.00007262380328602412
$ws.Cells.Item(13,16).Value = 0.00007262380328602412
$ws.Cells.Item(13,17).Value = 0.001154525204
$ws.Cells.Item(13,18).Value = 0.010390726836
$ws.Cells.Item(13,19).Value = 0.00001303719771837014
$ws.Cells.Item(13,20).Value = 0.00001303719771837014

$ws.Cells.Item(14,1).Value = 'Resolving-Mac'
$ws.Cells.Item(14,2).Value = 'Dnajb11'
$ws.Cells.Item(14,3).Value = 'Prtg'
$ws.Cells.Item(14,4).Value = 'ECs'
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 11.269395
$ws.Cells.Item(14,8).Value = 33.808185
$ws.Cells.Item(14,9).Value = 0.2821163927527693
$ws.Cells.Item(14,10).Value = 0.2821163927527693
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.1416426666666667
$ws.Cells.Item(14,14).Value = 0.424928
$ws.Cells.Item(14,15).Value = 0.06389210658948997
$ws.Cells.Item(14,16).Value = 0.06389210658948997
$ws.Cells.Item(14,17).Value = 1.59622715952
$ws.Cells.Item(14,18).Value = 14.36604443568
$ws.Cells.Item(14,19).Value = 0.01802501063640235
$ws.Cells.Item(14,20).Value = 0.01802501063640235

$ws.Cells.Item(15,1).Value = 'Resolving-Mac'
$ws.Cells.Item(15,2).Value = 'Dnajb11'
$ws.Cells.Item(15,3).Value = 'Prtg'
$ws.Cells.Item(15,4).Value = 'FAPs'
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 11.269395
$ws.Cells.Item(15,8).Value = 33.808185
$ws.Cells.Item(15,9).Value = 0.2821163927527693
$ws.Cells.Item(15,10).Value = 0.2821163927527693
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 1.357753333333333
$ws.Cells.Item(15,14).Value = 4.073259999999999
$ws.Cells.Item(15,15).Value = 0.6124547266518232
$ws.Cells.Item(15,16).Value = 0.6124547266518232
$ws.Cells.Item(15,17).Value = 15.3010586259
$ws.Cells.Item(15,18).Value = 137.7095276331
$ws.Cells.Item(15,19).Value = 0.1727835182073957
$ws.Cells.Item(15,20).Value = 0.1727835182073957

$ws.Cells.Item(16,1).Value = 'Resolving-Mac'
$ws.Cells.Item(16,2).Value = 'Dnajb11'
$ws.Cells.Item(16,3).Value = 'Prtg'
$ws.Cells.Item(16,4).Value = 'MuSCs'
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 11.269395
$ws.Cells.Item(16,8).Value = 33.808185
$ws.Cells.Item(16,9).Value = 0.2821163927527693
$ws.Cells.Item(16,10).Value = 0.2821163927527693
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.7173469999999998
$ws.Cells.Item(16,14).Value = 2.152041
$ws.Cells.Item(16,15).Value = 0.3235805429554008
$ws.Cells.Item(16,16).Value = 0.3235805429554009
$ws.Cells.Item(16,17).Value = 8.084066695064999
$ws.Cells.Item(16,18).Value = 72.75660025558498
$ws.Cells.Item(16,19).Value = 0.0912873755435602
$ws.Cells.Item(16,20).Value = 0.09128737554356021

$ws.Cells.Item(17,1).Value = 'Resolving-Mac'
$ws.Cells.Item(17,2).Value = 'Dnajb11'
$ws.Cells.Item(17,3).Value = 'Prtg'
$ws.Cells.Item(17,4).Value = 'Resolving-Mac'
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 11.269395
$ws.Cells.Item(17,8).Value = 33.808185
$ws.Cells.Item(17,9).Value = 0.2821163927527693
$ws.Cells.Item(17,10).Value = 0.2821163927527693
$ws.Cells.Item(17,11).Value = 1
$ws.Cells.Item(17,12).Value = 0.3333333333333333
$ws.Cells.Item(17,13).Value = 0.000161
$ws.Cells.Item(17,14).Value = 0.000483
$ws.Cells.Item(17,15).Value = 0.00007262380328602412
$ws.Cells.Item(17,16).Value = 0.00007262380328602412
$ws.Cells.Item(17,17).Value = 0.001814372595
$ws.Cells.Item(17,18).Value = 0.016329353355
$ws.Cells.Item(17,19).Value = 0.00002048836541103984
$ws.Cells.Item(17,20).Value = 0.00002048836541103984
